# FIX: Adds correct alignment (#89)
# Adds a new "alignment" worksheet after Sheet1 demonstrating the various
# horizontal/vertical alignment options, and makes it the active sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

$ws = $wb.Worksheets.Add([System.Type]::Missing, $ws1)
$ws.Name = "alignment"

# --- Column widths ----------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 19.998697916666668
$ws.Columns.Item(2).ColumnWidth = 15.498697916666666

# --- Row heights --------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 32
$ws.Rows.Item(3).RowHeight = 39
$ws.Rows.Item(4).RowHeight = 42
$ws.Rows.Item(5).RowHeight = 56
$ws.Rows.Item(6).RowHeight = 55
$ws.Rows.Item(7).RowHeight = 59
$ws.Rows.Item(8).RowHeight = 56

# --- Row 1/2: plain header + "no alignment" demo row -------------------
$ws.Range("A1").Value = "Horizontal"
$ws.Range("B1").Value = "Vertical"
$ws.Range("A1:B1").Font.Bold = $true

$ws.Range("A2").Value = "No alignment"
$ws.Range("B2").Value = "No alignment"

# --- Column A: horizontal alignment demos (rows 3-8) --------------------
$ws.Range("A3").Value = "Left"
$ws.Range("A3").HorizontalAlignment = -4131

$ws.Range("A4").Value = "Distributed"
$ws.Range("A4").HorizontalAlignment = -4117

$ws.Range("A5").Value = "Right"
$ws.Range("A5").HorizontalAlignment = -4152

$ws.Range("A6").Value = "Center"
$ws.Range("A6").HorizontalAlignment = -4108

$ws.Range("A7").Value = "Fill"
$ws.Range("A7").HorizontalAlignment = 5

$ws.Range("A8").Value = "Justify"
$ws.Range("A8").HorizontalAlignment = -4130

# --- Column B: vertical alignment demos (rows 3-7) -----------------------
$ws.Range("B3").Value = "Bottom"

$ws.Range("B4").Value = "Top"
$ws.Range("B4").VerticalAlignment = -4160

$ws.Range("B5").Value = "Center"
$ws.Range("B5").VerticalAlignment = -4108

$ws.Range("B6").Value = "Justify"
$ws.Range("B6").VerticalAlignment = -4130

$ws.Range("B7").Value = "Distributed"
$ws.Range("B7").VerticalAlignment = -4117

# --- Selection / activation ----------------------------------------------
$ws.Select() | Out-Null
$ws.Range("B6").Select() | Out-Null
